# Apportion share of building retrofitting and distributed solar revenues to labor (#19)

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename the "SYSoCCtaSC" sheet to "SYSoCCtaSC-electricity" and add a
#    new "SYSoCCtaSC-buildings" sheet right after it, mirroring its
#    layout for the buildings / distributed-solar retrofit variable.
#    (Done before the "About" edit below so new shared strings land in
#    the same order as the source edit: "distributed solar" first.)
# ------------------------------------------------------------------
$wsElec = $wb.Worksheets.Item("SYSoCCtaSC")
$wsElec.Name = "SYSoCCtaSC-electricity"

$wsBuild = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsElec)
$wsBuild.Name = "SYSoCCtaSC-buildings"
# Same tab color as the other SYSoCCtaSC sheet (theme Accent5, Darker 25%).
$wsBuild.Tab.ThemeColor = 8
$wsBuild.Tab.TintAndShade = -0.249977111117893

$wsBuild.Range("B1").Value = $wsElec.Range("B1").Text
$wsBuild.Range("B1").Style = $wsElec.Range("B1").Style

$wsBuild.Range("A2").Value = "distributed solar"
$wsBuild.Range("B2").Formula = "=AVERAGE(Data!B64:B65)"
$wsBuild.Range("B2").Style = $wsElec.Range("B2").Style

$wsBuild.Columns.Item(1).ColumnWidth = $wsElec.Columns.Item(1).ColumnWidth
$wsBuild.Columns.Item(2).ColumnWidth = $wsElec.Columns.Item(2).ColumnWidth

# ------------------------------------------------------------------
# 2. "About" sheet: insert a new row 2 describing the new variable,
#    pushing the existing content (and its hyperlinks) down by one row.
# ------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

# Remember current hyperlink anchors (row/col/address) before the shift,
# since row-insert does not renumber existing Hyperlink ranges.
$links = @()
foreach ($hl in $wsAbout.Hyperlinks) {
    $links += , @($hl.Range.Row, $hl.Range.Column, $hl.Address)
}
$wsAbout.Hyperlinks.Delete()

$wsAbout.Rows.Item(2).Insert()
$wsAbout.Range("A2").Value = "SYSoCCtaSC Share of Distributed Solar and Retrofitting Costs that is Labor"
$wsAbout.Range("A2").Font.Bold = $true

# Re-create the hyperlinks one row lower than before, then restore the
# original (non-"Hyperlink") cell style that .Hyperlinks.Add() overwrites.
foreach ($l in $links) {
    $r = [int]$l[0] + 1
    $c = [int]$l[1]
    $addr = $l[2]
    $cell = $wsAbout.Cells.Item($r, $c)
    $styleBefore = $cell.Style
    $wsAbout.Hyperlinks.Add($cell, $addr) | Out-Null
    $cell.Style = $styleBefore
}
